# The cross column previously lumped the F1 and F2 generations together
# (two rows each labelled "F1" / "F2"). Split them into distinct
# sub-categories: F1a/F1b and F2a/F2b. The old "P2" row (last row) keeps
# its label; only the underlying shared-string bookkeeping shifts.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "F1a"
$ws.Range("A4").Value = "F1b"
$ws.Range("A5").Value = "F2a"
$ws.Range("A6").Value = "F2b"
$ws.Range("A7").Value = "P2"

# Restore the on-save selection/view state (cursor on E12).
$ws.Activate()
$ws.Range("E12").Select()
